{"js": "// The author typed the word \"code&\" in front of \"views work\" inside the\n// sentence \"... enlighten us about how/whether the views work\".\n// That single text edit is what moved Word's automatic \"_GoBack\" bookmark\n// (which always tracks the most recent edit position) from its previous\n// location (after the two <w:br/> runs, right before \"*) \u201cXML-like\u201d...\")\n// to the new edit point (right after \"code&\", before \"views work\").\n\nconst body = context.document.body;\n\n// 1) Insert \"code&\" right before \"views work\" in the target sentence.\nconst target = body.search(\"views work\", { matchCase: true, matchWholeWord: false });\ntarget.load(\"text\");\nawait context.sync();\n\nif (target.items.length === 0) {\n  throw new Error('Could not find \"views work\" in the document body.');\n}\n\n// Use the first (and only) match.\nconst viewsWorkRange = target.items[0];\nviewsWorkRange.insertText(\"code&\", Word.InsertLocation.before);\nawait context.sync();\n\n// 2) The inserted word should live in its own run, distinct from the run\n// that held the untouched \"... the \" text before it, and from the run\n// that holds \"views work\" after it. Re-locate the merged text and\n// re-insert it as two explicit runs via OOXML so the \"code&\" word is its\n// own run (matching how Word materializes a freshly typed word).\nconst merged = body.search(\"code&views work\", { matchCase: true });\nmerged.load(\"text\");\nawait context.sync();\n\nif (merged.items.length > 0) {\n  const mergedRange = merged.items[0];\n  const splitOoxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:rPr><w:b/><w:i/><w:lang w:val=\"en-US\"/></w:rPr><w:t>code&amp;</w:t></w:r>' +\n    '<w:r><w:rPr><w:b/><w:i/><w:lang w:val=\"en-US\"/></w:rPr><w:t>views work</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n  mergedRange.insertOoxml(splitOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Move the \"_GoBack\" bookmark: delete it from wherever it currently is\n// (its old position, before \"*) \u201cXML-like\u201d...\") and re-insert it right\n// between the new \"code&\" run and the \"views work\" run.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst viewsWork2 = body.search(\"views work\", { matchCase: true });\nviewsWork2.load(\"text\");\nawait context.sync();\n\nif (viewsWork2.items.length > 0) {\n  const insertionPoint = viewsWork2.items[0].getRange(Word.RangeLocation.start);\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The author typed the word \"code&\" in front of \"views work\" inside the\n# sentence \"... enlighten us about how/whether the views work\".\n# That single text edit is what moved Word's automatic \"_GoBack\" bookmark\n# (which always tracks the most recent edit position) from its previous\n# location (after the two Br runs, right before \"*) \"XML-like\"...\")\n# to the new edit point (right after \"code&\", before \"views work\").\n\n$d = $word.ActiveDocument\n\n# 1) Insert \"code&\" right before \"views work\" in the target sentence.\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.Text = \"views work\"\n$found = $find.Execute()\nif (-not $found) {\n    throw 'Could not find \"views work\" in the document body.'\n}\n$searchRange.InsertBefore(\"code&\")\n\n# 2) Re-split the now-merged \"code&views work\" text into two explicit\n# runs - one for the freshly typed \"code&\" and one for the pre-existing\n# \"views work\" - so the new word lives in its own run, matching how Word\n# materializes a freshly typed word next to existing text.\n$splitRange = $d.Content\n$splitFind = $splitRange.Find\n$splitFind.Text = \"code&views work\"\n$splitFound = $splitFind.Execute()\nif ($splitFound) {\n    $target = $d.Range($splitRange.Start, $splitRange.End)\n    $ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:rPr><w:b/><w:i/><w:lang w:val=\"en-US\"/></w:rPr><w:t>code&amp;</w:t></w:r><w:r><w:rPr><w:b/><w:i/><w:lang w:val=\"en-US\"/></w:rPr><w:t>views work</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $target.InsertXML($ooxml)\n}\n\n# 3) Move the \"_GoBack\" bookmark: delete it from wherever it currently is\n# (its old position, before \"*) \"XML-like\"...\") and re-insert it right\n# between the new \"code&\" run and the \"views work\" run.\n$bookmarks = $d.Bookmarks\ntry {\n    $existing = $bookmarks.Item(\"_GoBack\")\n    $existing.Delete()\n} catch {\n    # No existing _GoBack bookmark - nothing to remove.\n}\n\n$finalRange = $d.Content\n$finalFind = $finalRange.Find\n$finalFind.Text = \"views work\"\n$finalFound = $finalFind.Execute()\nif ($finalFound) {\n    $insertionPoint = $d.Range($finalRange.Start, $finalRange.Start)\n    $bookmarks.Add(\"_GoBack\", $insertionPoint)\n}\n"}
